$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume(1h) (E) columns for the cells we touch
# so that Excel does not reinterpret numeric-looking strings as numbers,
# and percent-looking strings as percentages.


# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.896.15'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.412.07'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.19%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.51'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.55%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.50'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +8.44%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.411.07'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.13%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.15%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.49'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.14%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +9.50%  '

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.32%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.998.22'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.61%  '

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +8.06%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.414.30'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.74%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.46'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.39%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.955.24'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.91%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.12'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +6.82%  '

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.00%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.77%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.96'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +11.21%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.571'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.46%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.550.19'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.39%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000129'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +17.61%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.05%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.48'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.38%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.51%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +8.21%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.53%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.30'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.69%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.71%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.56%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.443.09'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.27%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.16%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.46'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.06%  '

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.77%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.82%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '163.39'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0791'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.10%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +15.62%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +7.07%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.21%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.98%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.47'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.42%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.11'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +9.71%  '

# Row 48
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '41.78'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.02%  '

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.05%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.09'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +6.21%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.381.54'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +11.20%  '
